$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "BZT"
$ws.Range("B15").Value = "benzetacil "
$ws.Range("C15").Value = 19
$ws.Range("D15").Value = 12
